$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'DExpi910'
$ws.Range("B2").Value = 231006242
$ws.Range("C2").Value = 'xsxsafl56'
$ws.Range("D2").Value = 'f7R9F!%y'
$ws.Range("E2").Value = 'MR'
$ws.Range("F2").Value = 'TEevobFl'
$ws.Range("G2").Value = 'BUbp'
$ws.Range("H2").Value = 'Candidate'
$ws.Range("A3").Value = 'dEtqV203'
$ws.Range("B3").Value = 231006241
$ws.Range("C3").Value = 'qlbigcc48'
$ws.Range("D3").Value = 'a5GQd6&#'
$ws.Range("E3").Value = 'MR'
$ws.Range("F3").Value = 'LZSOyssD'
$ws.Range("G3").Value = 'XKkU'
$ws.Range("H3").Value = 'Candidate'
$ws.Range("A4").Value = 'eugBU870'
$ws.Range("B4").Value = 231006240
$ws.Range("C4").Value = 'jtrxryw40'
$ws.Range("D4").Value = 'r$37eTS#'
$ws.Range("E4").Value = 'MR'
$ws.Range("F4").Value = 'QcgKJLVH'
$ws.Range("G4").Value = 'GYan'
$ws.Range("H4").Value = 'Candidate'
$ws.Range("A5").Value = 'fSaxA904'
$ws.Range("B5").Value = 231006239
$ws.Range("C5").Value = 'nydxdua41'
$ws.Range("D5").Value = 'sE5d!G3#'
$ws.Range("E5").Value = 'MR'
$ws.Range("F5").Value = 'ZRTBOawx'
$ws.Range("G5").Value = 'QNWF'
$ws.Range("H5").Value = 'Candidate'
$ws.Range("A6").Value = 'eBFVJ660'
$ws.Range("B6").Value = 231006238
$ws.Range("C6").Value = 'japvssm93'
$ws.Range("D6").Value = 'JzS%68#h'
$ws.Range("E6").Value = 'MR'
$ws.Range("F6").Value = 'IhkLgPlU'
$ws.Range("G6").Value = 'GdnO'
$ws.Range("H6").Value = 'Candidate'
$ws.Range("A7").Value = 'dVPSh786'
$ws.Range("B7").Value = 231006237
$ws.Range("C7").Value = 'fnoefug22'
$ws.Range("D7").Value = 'py7&Z8A%'
$ws.Range("E7").Value = 'MR'
$ws.Range("F7").Value = 'lHDyDGKT'
$ws.Range("G7").Value = 'oIvn'
$ws.Range("H7").Value = 'Candidate'
$ws.Range("A8").Value = 'ONeky198'
$ws.Range("B8").Value = 231006236
$ws.Range("C8").Value = 'wxqmymu93'
$ws.Range("D8").Value = 'Ek6F5#%t'
$ws.Range("E8").Value = 'MR'
$ws.Range("F8").Value = 'SPjvPwBK'
$ws.Range("G8").Value = 'LOrG'
$ws.Range("H8").Value = 'Candidate'
$ws.Range("A9").Value = 'dUUGn232'
$ws.Range("B9").Value = 231006235
$ws.Range("C9").Value = 'cbhryaz47'
$ws.Range("D9").Value = 'YdS!42&y'
$ws.Range("E9").Value = 'MR'
$ws.Range("F9").Value = 'ZrLZtrhj'
$ws.Range("G9").Value = 'MQQV'
$ws.Range("H9").Value = 'Candidate'
$ws.Range("A10").Value = 'yboiX674'
$ws.Range("B10").Value = 231006234
$ws.Range("C10").Value = 'ekrhkwb77'
$ws.Range("D10").Value = 'TSzf#49&'
$ws.Range("E10").Value = 'MR'
$ws.Range("F10").Value = 'gIYiuzuk'
$ws.Range("G10").Value = 'JDHp'
$ws.Range("H10").Value = 'Candidate'
$ws.Range("A11").Value = 'mtPcK694'
$ws.Range("B11").Value = 231006233
$ws.Range("C11").Value = 'iynpfep44'
$ws.Range("D11").Value = 'y#&P3K5g'
$ws.Range("E11").Value = 'MR'
$ws.Range("F11").Value = 'oBeZMazR'
$ws.Range("G11").Value = 'NCKN'
$ws.Range("H11").Value = 'Candidate'
$ws.Range("A12").Value = 'mjriQ397'
$ws.Range("B12").Value = 231006232
$ws.Range("C12").Value = 'zzuumhx76'
$ws.Range("D12").Value = 'w%7&6FYs'
$ws.Range("E12").Value = 'MR'
$ws.Range("F12").Value = 'uCgIhQZi'
$ws.Range("G12").Value = 'faau'
$ws.Range("H12").Value = 'Candidate'
$ws.Range("A13").Value = 'jDuJK351'
$ws.Range("B13").Value = 231006231
$ws.Range("C13").Value = 'zvdeuav30'
$ws.Range("D13").Value = 'gp&%H7B4'
$ws.Range("E13").Value = 'MR'
$ws.Range("F13").Value = 'jNSTiFse'
$ws.Range("G13").Value = 'JLww'
$ws.Range("H13").Value = 'Candidate'

$ws.Range("A8:H13").Borders.LineStyle = 1
$ws.Range("A1:H13").Select()
